$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = "2020-09-30 00:00:00"

$ws.Range("O2").Value = 49132470.49
$ws.Range("P2").Value = 247689138.13
$ws.Range("Q2").Value = 191999725.34
$ws.Range("R2").Value = -11.1076778224
$ws.Range("S2").Value = 153094704.74
$ws.Range("T2").Value = 153094704.74
$ws.Range("U2").Value = -9.0215776972
$ws.Range("V2").Value = 12033616.75
$ws.Range("W2").Value = 16079754.47
$ws.Range("X2").Value = 346855.67
$ws.Range("Y2").Value = 57079843.78
$ws.Range("Z2").Value = 56806982.27
$ws.Range("AA2").Value = 7674511.78
$ws.Range("AG2").Value = 2858424.12
$ws.Range("AP2").Value = -6.6950813206
$ws.Range("AQ2").Value = 13.540271063682
$ws.Range("AR2").Value = 16.828594620139
$ws.Range("AS2").Value = 45354963.11
$ws.Range("AT2").Value = 9.174594331711999
